$filesQuery = @"
MATCH (f:file)-->(s:study)
OPTIONAL MATCH (samp:sample)<--(f)
OPTIONAL MATCH (samp)-->(p:participant)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Targeted Sequencing" IN es
WITH DISTINCT f, s, p, samp
RETURN
    coalesce(f.file_name, '') as ``File Name``,
    coalesce(s.study_name,'') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(p.participant_id, '') as ``Participant ID``,
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(f.file_type, '') as ``File Type``
ORDER BY f.file_name LIMIT 100
"@

$samplesQuery = @"
MATCH (samp:sample)-->(p:participant)-->(s:study)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
OPTIONAL MATCH (p)<--(diag:diagnosis)
WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Targeted Sequencing" IN es
WITH DISTINCT s, p, samp
RETURN
    coalesce(samp.sample_id, '') as ``Sample ID``,
    coalesce(p.participant_id,'') as ``Participant ID``,
    coalesce(s.study_name, '') as ``Study Name``,
    coalesce(s.phs_accession,'') as ``Accession``,
    coalesce(samp.sample_tumor_status,'') as ``Tumor``,
    coalesce(samp.sample_type,'') as ``Analyte Type``
ORDER BY samp.sample_id LIMIT 100
"@

$participantsQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (samp)<--(f:file)
WITH p, samp, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
WHERE "Targeted Sequencing" IN es
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct coalesce(samp.sample_id, "Not specified in data"))) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id LIMIT 100
"@

$statQuery = @"
CALL{
    MATCH (p:participant)-->(s:study)
    OPTIONAL MATCH (samp:sample)-->(p)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Targeted Sequencing" IN es
    RETURN 
        count(distinct p) AS num_participants
}
WITH num_participants
CALL {
    MATCH (samp:sample)-->(p:participant)-->(s:study)
    OPTIONAL MATCH (samp)<--(f:file)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Targeted Sequencing" IN es
    RETURN 
        count(distinct samp) AS num_samples
}
WITH num_participants, num_samples
CALL {
    MATCH (f:file)-->(s:study)
    OPTIONAL MATCH (f)<--(g:genomic_info)
    OPTIONAL MATCH (samp:sample)<--(f)
    OPTIONAL MATCH (p:participant)<--(samp)
    OPTIONAL MATCH (p)<--(diag:diagnosis)
    WITH s, p, samp, f, g, diag, apoc.coll.flatten(COLLECT (apoc.text.split(f.experimental_strategy_and_data_subtypes,"[;,]\s{0,1}")), true) as es
    WHERE "Targeted Sequencing" IN es
    RETURN 
        count(distinct s) AS num_studies,
        count(distinct f) AS num_files
}
RETURN 
    num_studies AS Studies,
    num_participants AS Participants,
    num_samples AS Samples,
    num_files AS ``Files``
"@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update query text in the three tab rows ---
# NOTE: new shared-strings are appended in the order the .Value assignments
# are made (first-seen order), so we deliberately write Files, then Samples,
# then Participants, then StatQuery to reproduce the target shared string order.
$ws.Range("B4").Value = $filesQuery
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B2").Value = $participantsQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("C3").Value = $statQuery
$ws.Range("C4").Value = $statQuery

# --- Row heights (wrap-text cells need explicit custom heights) ---
$ws.Rows(2).RowHeight = 328.5
$ws.Rows(3).RowHeight = 294.75
$ws.Rows(4).RowHeight = 297.75

# --- Update sheet view: scroll position + selection ---
$ws.Range("A3").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E14").Select()
